$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update row 28 styling: switch from style 1 to style 5 (fontId1,fillId0,borderId4) ---
# B5 already carries style 5 in the workbook; copy only its format onto the whole row 28 range.
$fmtSrc5 = $ws.Range("B5")
$row28 = $ws.Range("A28:D28")
$fmtSrc5.Copy()
$row28.PasteSpecial(-4122)

# --- 2. Add the 8 new rows (29-36) with new beneficial bacteria entries ---
# New NCBI/MIrROR name pairs (shared strings 53-68) and the '유익' (beneficial) sign in column D.
$newRows = @(
    @{ Row = 29; NCBI = "Lactobacillus acidophilus";    Mirror = "s__Lactobacillus_acidophilus";    Height = 38 },
    @{ Row = 30; NCBI = "Lactobacillus plantarum";       Mirror = "s__Lactobacillus_plantarum";       Height = 38 },
    @{ Row = 31; NCBI = "Bifidobacterium animalis";      Mirror = "s__Bifidobacterium_animalis";      Height = 38 },
    @{ Row = 32; NCBI = "Lactobacillus reuteri";         Mirror = "s__Lactobacillus_reuteri";         Height = 38 },
    @{ Row = 33; NCBI = "Lactobacillus casei";           Mirror = "s__Lactobacillus_casei";           Height = 38 },
    @{ Row = 34; NCBI = "Bifidobacterium breve";         Mirror = "s__Bifidobacterium_breve";         Height = 38 },
    @{ Row = 35; NCBI = "Bifidobacterium bifidum";       Mirror = "s__Bifidobacterium_bifidum";       Height = 38 },
    @{ Row = 36; NCBI = "Streptococcus thermophilus";    Mirror = "s__Streptococcus_thermophilus";    Height = 50.5 }
)

# Template row (27) already uses plain style 1 across all four columns A:D - reuse its formatting.
$fmtSrc1 = $ws.Range("A27:D27")

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Apply the standard formatting (style 1) to the whole row first.
    $dst = $ws.Range("A" + $r + ":D" + $r)
    $fmtSrc1.Copy()
    $dst.PasteSpecial(-4122)

    # Fill in the cell values.
    $ws.Cells.Item($r, 1).Value = $entry.NCBI
    $ws.Cells.Item($r, 2).Value = $entry.Mirror
    $ws.Cells.Item($r, 4).Value = "유익"

    # Match the row height used in the source layout.
    $ws.Rows.Item($r).RowHeight = $entry.Height
}

# Row 31 (Bifidobacterium animalis) highlights the NCBI/MIrROR columns with a white fill (new style).
$ws.Range("A31:B31").Interior.Color = 16777215

# --- 3. Update the active selection to match the saved view ---
$ws.Range("G3").Select()
